# Refresh the crypto price/volume table with the latest scraped values.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h).
# Note: several Price cells look like plain numbers (e.g. "1.00", "6.30").
# Those columns are plain text in the workbook (no thousands grouping, and
# some prices use "." as a thousands separator, e.g. "67.879.12"), so for
# cells whose new value would otherwise be auto-parsed as a number by
# Excel we force NumberFormat "@" (Text) first to keep them as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.879.12'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '3.850.99'
$ws.Range("E3").Value = '  -1.49%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.58'
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.34'
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("D7").Value = '3.850.26'
$ws.Range("E7").Value = '  -1.50%  '
$ws.Range("E9").Value = '  -0.72%  '
$ws.Range("E10").Value = '  -0.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.30'
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("E13").Value = '  +1.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.91'
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").Value = '4.500.48'
$ws.Range("E15").Value = '  -1.34%  '
$ws.Range("D16").Value = '3.839.14'
$ws.Range("E16").Value = '  -2.15%  '
$ws.Range("D17").Value = '67.963.45'
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.19'
$ws.Range("E18").Value = '  +7.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.42'
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("E20").Value = '  -1.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.79'
$ws.Range("E21").Value = '  -3.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '466.84'
$ws.Range("E22").Value = '  -3.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.729'
$ws.Range("E23").Value = '  +1.32%  '
$ws.Range("E24").Value = '  -4.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.19'
$ws.Range("E25").Value = '  -1.42%  '
$ws.Range("E26").Value = '  -1.30%  '
$ws.Range("E27").Value = '  +0.89%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.94'
$ws.Range("E30").Value = '  +0.28%  '
$ws.Range("D31").Value = '4.003.86'
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.67'
$ws.Range("E32").Value = '  -2.53%  '
$ws.Range("E33").Value = '  -3.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.94'
$ws.Range("E34").Value = '  -3.09%  '
$ws.Range("D35").Value = '3.832.22'
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("E36").Value = '  -2.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.140'
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("E38").Value = '  -2.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.89'
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.25'
$ws.Range("E40").Value = '  +8.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("E42").Value = '  -2.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '426.43'
$ws.Range("E43").Value = '  -3.00%  '
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.14'
$ws.Range("E46").Value = '  -2.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.53'
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '143.68'
$ws.Range("E48").Value = '  +1.14%  '
$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000270'
$ws.Range("E49").Value = '  +11.85%  '
$ws.Range("B50").Value = 'Arweave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '39.43'
$ws.Range("E50").Value = '  +0.67%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.80'
$ws.Range("E51").Value = '  -0.62%  '
